# Apply the "updates to model files" commit:
#  - L2/M2 (Polymyxins/Other price inputs) get corrected values
#  - K4/K5 (Pleuromutilins group inputs) get corrected values
#  - H8 keeps the same Aminoglycosides revenue formula (H2*H6); the ripple
#    from L2/M2 recomputes L8/M8 automatically via the existing formulas
#  - selection moves to G20

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("L2").Value = 980.58
$ws.Range("M2").Value = 5984.73

$ws.Range("K4").Value = 2
$ws.Range("K5").Value = 2

$ws.Range("H8").Formula = "=H2*H6"

$ws.Range("G20").Select()
